$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New label row under the existing table (mirrors D2's "Time taken in hours:")
$ws.Range("D9").Value = "Time taken in hours:"

# New UCT sqrt2 vs UCT 0.5 result row
$ws.Range("B12").Value = "MCTS 10k Sarasua 1 UCT sqrt2 vs MCTS 10k Sarasua 1 UCT 0.5"
$ws.Range("C12").Value = "10.5/20"
$ws.Range("D12").Formula = "=8059/60/60"

# New UCT sqrt2 vs UCT 5 result row
$ws.Range("B13").Value = "MCTS 10k Sarasua 1 UCT sqrt2 vs MCTS 10k Sarasua 1 UCT 5"
$ws.Range("C13").Value = "11.0/20"
$ws.Range("C13").NumberFormat = "d-mmm"
$ws.Range("D13").Formula = "=10090/3600"

# Column B needs to widen to fit the new long labels (bestFit-style autosize
# to the new longer strings added below)
$ws.Columns.Item(2).ColumnWidth = 54.67

# Match the author's final selection
$ws.Range("B11").Select()
